# Update the "想去人数" (interested count) values on the "展览" and
# "全部类型" worksheets, matching the generated-output refresh captured
# in the commit.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 11
    $ws.Range("F4").Value = 954
}
